# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
# Swap the match-record rows that had been mismatched: the full data
# (columns B through AB; column A is just the running index and is left
# untouched) of row 37 <-> row 38, and of row 148 <-> row 149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($sheet, $rowA, $rowB, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $sheet.Cells.Item($rowA, $col)
        $cellB = $sheet.Cells.Item($rowB, $col)

        $valA = $cellA.Value()
        $valB = $cellB.Value()

        if ($valA -ne $valB) {
            $cellA.Value = $valB
            $cellB.Value = $valA
        }
    }
}

# Column B = 2 ... Column AB = 28
Swap-RowData $ws 37 38 2 28
Swap-RowData $ws 148 149 2 28
